# Regenerate s_val data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-8.
# Column A (date) and F (Win) remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243; E = 0.4998867070740569; G = 22.31973251085698 }
    3 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243; E = 0.4998867070740569; G = 22.31973251085698 }
    5 = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 157.8057217802531; E = 6.48142807727062;   G = 164.7801781670697 }
    6 = @{ B = 0.7287194209349384; C = 0.3375848360084654; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.719096746035642 }
    7 = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 3.082599426703578; E = 0.4998867070740569; G = 5.964442013611383 }
    8 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 3.082599426703578; E = 0.4998867070740569; G = 6.741336633845642 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}

$wb.Save()
